# New PO forecast model
# Updates three sheets: "Weekly Quantity", "Monthly Trend", "PO Forecast".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append two new weekly rows (26, 27)
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @(45662.99999999999, 5),
    @(45669.99999999999, 1)
)

$r = 26
foreach ($row in $weeklyNewRows) {
    $wsWeekly.Cells.Item($r, 1).Value = $row[0]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = $wsWeekly.Cells.Item($r - 1, 1).NumberFormat
    $wsWeekly.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append one new monthly row (16)
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Cells.Item(16, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(16, 1).NumberFormat = $wsMonthly.Cells.Item(15, 1).NumberFormat
$wsMonthly.Cells.Item(16, 2).Value = 6

# ---------------------------------------------------------------------------
# Sheet 3: "PO Forecast" - new forecast model: revised values for existing
# rows 2-33 plus two additional forecast rows (34, 35)
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastRows = @(
    @(44983.99999999999, 5),
    @(44997.99999999999, 5),
    @(45004.99999999999, 5),
    @(45011.99999999999, 5),
    @(45039.99999999999, 5),
    @(45088.99999999999, 5),
    @(45095.99999999999, 5),
    @(45102.99999999999, 5),
    @(45193.99999999999, 5),
    @(45200.99999999999, 5),
    @(45207.99999999999, 5),
    @(45221.99999999999, 5),
    @(45235.99999999999, 5),
    @(45452.99999999999, 6),
    @(45466.99999999999, 6),
    @(45487.99999999999, 6),
    @(45508.99999999999, 6),
    @(45515.99999999999, 6),
    @(45543.99999999999, 6),
    @(45578.99999999999, 6),
    @(45585.99999999999, 6),
    @(45606.99999999999, 6),
    @(45634.99999999999, 6),
    @(45641.99999999999, 6),
    @(45662.99999999999, 6),
    @(45669.99999999999, 7),
    @(45676.99999999999, 7),
    @(45683.99999999999, 7),
    @(45690.99999999999, 7),
    @(45697.99999999999, 7),
    @(45704.99999999999, 7),
    @(45711.99999999999, 7),
    @(45718.99999999999, 7),
    @(45725.99999999999, 7)
)

$r = 2
foreach ($row in $forecastRows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
